{"js": "// The author placed the cursor between \"M.\" and \"ammadov\" in the\n// \"Tilakaratne, C., M. Mammadov, ...\" reference entry and saved the\n// document. Word stamps that last-edit location with its automatic\n// \"_GoBack\" bookmark (used for Shift+F5 / \"return to last edit\").\n// Reproduce that by locating the text and dropping a collapsed\n// bookmark named \"_GoBack\" right after \"Tilakaratne, C., M. M\".\n\nconst searchResults = context.document.body.search(\"Tilakaratne, C., M. M\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const found = searchResults.items[0];\n  const endRange = found.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author placed the cursor between \"M.\" and \"ammadov\" in the\n# \"Tilakaratne, C., M. Mammadov, ...\" reference entry and saved the\n# document. Word stamps that last-edit location with its automatic\n# \"_GoBack\" bookmark (used for Shift+F5 / \"return to last edit\").\n# Reproduce that by locating the text and dropping a collapsed\n# bookmark named \"_GoBack\" right after \"Tilakaratne, C., M. M\".\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Tilakaratne, C., M. M\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $rng.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $rng)\n}\n"}
